$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column for "scenario" (mostly blank) and a new
# "headers" column right after "url" (old column B, new column C).
$ws.Columns("A").Insert()
$ws.Columns("D").Insert()

# Give the new "scenario" column (A) the same row-by-row formatting as
# the "method" column (B) that sits right next to it.
$ws.Range("B1:B13").Copy()
$ws.Range("A1:A13").PasteSpecial(-4122)

$ws.Range("A1").Value = "scenario"
$ws.Range("D1").Value = "headers"

# The POST /users and PUT /users/1 rows get a borderless placeholder in
# the new "headers" column (no header file needed for those examples).
$ws.Range("D5").Borders.LineStyle = -4142
$ws.Range("D5").IndentLevel = 0
$ws.Range("D7").Borders.LineStyle = -4142
$ws.Range("D7").IndentLevel = 0

# The /wallets row never got a placeholder cell in the new column.
$ws.Range("D13").Clear()

# New example row: GET /fx_rates/eur/usd with a request-headers file.
$ws.Range("A12:H12").Copy()
$ws.Range("A14:H14").PasteSpecial(-4122)
$ws.Range("B14").Value = "GET"
$ws.Range("C14").Value = "/fx_rates/eur/usd"
$ws.Range("D14").Value = "./fx_rates/request_headers.json"
$ws.Range("E14").Value = "{}"
$ws.Range("F14").Value = "{}"
$ws.Range("G14").Value = 200
$ws.Range("H14").Value = '{"rate": 1.12}'
